$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column G: "CurrencyCol" ---
$ws.Range("G1").Value = "CurrencyCol"
$ws.Columns.Item(7).ColumnWidth = 11.25

$currencyFmt = "[$$]* #,##0.00;[$$]* ""-""#,##0.00"" "";[$$]* ""-""??;@"

$ws.Range("G2").Value = 23
$ws.Range("G2").NumberFormat = $currencyFmt

$ws.Range("G3").Value = 1000
$ws.Range("G3").NumberFormat = $currencyFmt

$ws.Range("G4").Value = 100
$ws.Range("G4").NumberFormat = $currencyFmt

$ws.Range("G5").Value = 200
$ws.Range("G5").NumberFormat = $currencyFmt

$ws.Range("G6").Value = 230
$ws.Range("G6").NumberFormat = $currencyFmt

$ws.Range("G7").Value = 452
$ws.Range("G7").NumberFormat = $currencyFmt

$ws.Range("G8").Value = 573
$ws.Range("G8").NumberFormat = $currencyFmt

$ws.Range("G9").Value = 23
$ws.Range("G9").NumberFormat = $currencyFmt

$ws.Range("G10").Value = 2.45
$ws.Range("G10").NumberFormat = $currencyFmt

$ws.Range("G11").Value = 1021.56
$ws.Range("G11").NumberFormat = $currencyFmt

$ws.Range("G12").Value = 202.5
$ws.Range("G12").NumberFormat = $currencyFmt

# --- Row 8: drop the stray empty B8 cell ---
$ws.Range("B8").ClearContents()

# --- Row 13: D13 gets a very long decimal display format, G13 gets currency (builtin #7) ---
$ws.Range("D13").NumberFormat = "0.000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000"

$ws.Range("G13").Value = 33.223
$ws.Range("G13").NumberFormat = "$#,##0.00_);($#,##0.00)"

# --- New row 14 ---
$ws.Range("F14").Value = $false
$ws.Range("G14").Value = 22
$ws.Range("G14").NumberFormat = "$#,##0.00_);($#,##0.00)"

# --- Selection moves to G2 ---
$ws.Range("G2").Select() | Out-Null
